# 3.13 3rd Commit: update the email sending logic
# Append the newly-processed users (rows 4-6) to the "users" sheet, mirroring
# the rows that already exist on the "temp" sheet for Yijia, Jiangnan/JiangNan
# and Bellamy (aka "Art1st").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users")

# Helper: give a row's index cell (column A) the same bold / bordered /
# centered look already used by the existing index cells (A2, A3).
function Set-IndexCellStyle($cell) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1         # xlContinuous
    $cell.Borders.Weight = 2            # xlThin
}

# Row 4: Yijia Sun
$ws.Cells.Item(4, 1).Value = 2
Set-IndexCellStyle $ws.Cells.Item(4, 1)
$ws.Cells.Item(4, 2).Value = "Yijia"
$ws.Cells.Item(4, 3).Value = "Sun"
$ws.Cells.Item(4, 4).Value = "yijiasun@qq.com"
$ws.Cells.Item(4, 5).Value = "Shanghai"
$ws.Cells.Item(4, 6).Value = 8
$ws.Cells.Item(4, 7).Value = "Null"

# Row 5: JiangNan Cai
$ws.Cells.Item(5, 1).Value = 3
Set-IndexCellStyle $ws.Cells.Item(5, 1)
$ws.Cells.Item(5, 2).Value = "JiangNan"
$ws.Cells.Item(5, 3).Value = "Cai"
$ws.Cells.Item(5, 4).Value = "jiangnantsai404@gmail.com"
$ws.Cells.Item(5, 5).Value = "Los Angeles"
$ws.Cells.Item(5, 6).Value = -7
$ws.Cells.Item(5, 7).Value = "Null"

# Row 6: Art1st (Bellamy) Sun
$ws.Cells.Item(6, 1).Value = 4
Set-IndexCellStyle $ws.Cells.Item(6, 1)
$ws.Cells.Item(6, 2).Value = "Art1st"
$ws.Cells.Item(6, 3).Value = "Sun"
$ws.Cells.Item(6, 4).Value = "bellamy93158@gmail.com"
$ws.Cells.Item(6, 5).Value = "Los Angeles"
$ws.Cells.Item(6, 6).Value = -7
$ws.Cells.Item(6, 7).Value = "Null"

$wb.Save()
